$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks up front - their cell references don't
# track a later column insert, so we rebuild them from scratch once the
# final layout is in place.
$ws.Cells.Hyperlinks.Delete()

# Insert a new first column for "loginname" (shifts firstname/lastname/
# email/password one column to the right: B/C/D/E).
$ws.Range("A1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "loginname"
$ws.Range("B1").Value = "firstname"
$ws.Range("C1").Value = "lastname"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "password"
$ws.Range("F1").Value = "country"

# Row 2 - Romeo / Juliet
$ws.Range("A2").Value = "romeo@yahoo.com"
$ws.Range("B2").Value = "Romeo"
$ws.Range("C2").Value = "Juliet"
$ws.Range("D2").Value = "shakespeare@gmail.com"
$ws.Range("E2").Value = "Younglove"
$ws.Range("F2").Value = "United Kingdom"

# Row 3 - Dorian Gray
$ws.Range("A3").Value = "creapypainting@gmail.com"
$ws.Range("B3").Value = "Dorian"
$ws.Range("C3").Value = "Gray"
$ws.Range("D3").Value = "creapypainting@gmail.com"
$ws.Range("E3").Value = "Foreveryoung"
$ws.Range("F3").Value = "Ireland"

# Row 4 - Alesha Karamazov
$ws.Range("A4").Value = "goodness@gmail.com"
$ws.Range("B4").Value = "Alesha"
$ws.Range("C4").Value = "Karamazov"
$ws.Range("D4").Value = "goodness@gmail.com"
$ws.Range("E4").Value = "WarWorld"
$ws.Range("F4").Value = "Russia"

# Re-create the mailto hyperlinks on the login/email columns, in the same
# order as the target workbook (A2, D2, A3, D3, A4, D4).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:romeo@yahoo.com", "", "", "romeo@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:shakespeare@gmail.com", "", "", "shakespeare@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:creapypainting@gmail.com", "", "", "creapypainting@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:creapypainting@gmail.com", "", "", "creapypainting@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:goodness@gmail.com", "", "", "goodness@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:goodness@gmail.com", "", "", "goodness@gmail.com")

# Adding hyperlinks via COM auto-creates & applies the built-in "Hyperlink"
# named cell style (underline + theme color). Drop that named style and
# restore the original blue, non-underlined look (the same look the sheet
# already used for its one pre-existing hyperlink column) by repainting
# the format from that original hyperlink cell.
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("D2").Copy()
$ws.Range("A2:A4,D3:D4").PasteSpecial(-4122)
$ws.Range("A2").Value = "romeo@yahoo.com"
$ws.Range("A3").Value = "creapypainting@gmail.com"
$ws.Range("A4").Value = "goodness@gmail.com"
$ws.Range("D3").Value = "creapypainting@gmail.com"
$ws.Range("D4").Value = "goodness@gmail.com"

$ws.Range("A2").Select()
$excel.CutCopyMode = $false
